# UPCCodes.xlsx update
# - Renamed several item descriptions (shortened / restyled "Sauce:" prefixes, etc.)
# - Removed the old scratch row ("Cheese Tortelloni 16oz." / "Pumpkin Tortelloni 16oz." shifted
#   up by one row, "asdf"/123 test row removed)
# - Added a new "Sauce: Sage Butter" item
# - Added a block of " Ravioli: X" duplicate-name rows (for CTRL+M matching) reusing the
#   existing ravioli UPCs
# - Appended one more blank row at the bottom of the sheet
# - Moved the active selection to C32 (in progress of scrolling consistently)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Simple in-place text renames (UPC / B column untouched)
# ---------------------------------------------------------------------------
$renames = @{
    "A33" = "Jumbo Eggplant Ravioli"
    "A36" = "Cheese Ravioli"
    "A38" = "Lobster Ravioli"
    "A39" = "Meat + Spinach Ravioli"
    "A40" = "Mushroom Ravioli"
    "A41" = "Pumpkin Ravioli"
    "A44" = "Rice Pudding"
    "A48" = "Meat Lasagne"
    "A51" = "Sauce: Arrabiata"
    "A52" = "Sauce: Bolognese"
    "A53" = "Sauce: Pesto"
    "A54" = "Sauce: Marinara"
    "A55" = "Sauce: Mushroom Pesto"
    "A56" = "Sauce: Pink Sauce"
    "A57" = "Sauce: Puttanesca"
    "A58" = "Sauce: Sun-dried Tomato Pesto"
    "A59" = "Sauce: Tomato Basil"
}
foreach ($addr in $renames.Keys) {
    $ws.Range($addr).Value = $renames[$addr]
}

# ---------------------------------------------------------------------------
# 2. New row 60: "Sauce: Sage Butter" / 111 (copy formatting from row 59, the
#    last of the sauce rows, which already carries the item/UPC styles)
# ---------------------------------------------------------------------------
$ws.Range("A59:B59").Copy()
$ws.Range("A60:B60").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A60").Value = "Sauce: Sage Butter"
$ws.Range("B60").Value = 111

# ---------------------------------------------------------------------------
# 3. Old row 61 ("Cheese Tortelloni 16oz." / 692159005601) becomes a blank
#    separator row (format like row 34, the blank row above the ravioli block)
# ---------------------------------------------------------------------------
$ws.Range("A34:B34").Copy()
$ws.Range("A61:B61").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A61:B61").ClearContents()

# ---------------------------------------------------------------------------
# 4. Tortelloni rows shift down one, and drop their "16oz." suffix
# ---------------------------------------------------------------------------
$ws.Range("A62").Value = "Cheese Tortelloni"
$ws.Range("B62").Value = "692159005601"
$ws.Range("A63").Value = "Pumpkin Tortelloni"
$ws.Range("B63").Value = "692159005618"

# ---------------------------------------------------------------------------
# 5. Old row 64 ("asdf" / 123 scratch row) becomes a plain blank row (format
#    like row 73, matching the long run of blank rows further down)
# ---------------------------------------------------------------------------
$ws.Range("A73:B73").Copy()
$ws.Range("A64:B64").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A64:B64").ClearContents()

# ---------------------------------------------------------------------------
# 6. New " Ravioli: X" rows 65-72, reusing the UPCs already used by the
#    corresponding "X Ravioli" rows 35-42 (copy their formatting first)
# ---------------------------------------------------------------------------
$ws.Range("A35:B42").Copy()
$ws.Range("A65:B72").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A65").Value = " Ravioli: Cheese + Spinach"
$ws.Range("A66").Value = " Ravioli: Cheese"
$ws.Range("A67").Value = " Ravioli: Goat Cheese"
$ws.Range("A68").Value = " Ravioli: Lobster"
$ws.Range("A69").Value = " Ravioli: Meat + Spinach"
$ws.Range("A70").Value = " Ravioli: Mushroom"
$ws.Range("A71").Value = " Ravioli: Pumpkin"
$ws.Range("A72").Value = " Ravioli: Roasted Red Pepper"

$ws.Range("B65").Value = $ws.Range("B35").Value()
$ws.Range("B66").Value = $ws.Range("B36").Value()
$ws.Range("B67").Value = $ws.Range("B37").Value()
$ws.Range("B68").Value = $ws.Range("B38").Value()
$ws.Range("B69").Value = $ws.Range("B39").Value()
$ws.Range("B70").Value = $ws.Range("B40").Value()
$ws.Range("B71").Value = $ws.Range("B41").Value()
$ws.Range("B72").Value = $ws.Range("B42").Value()

# ---------------------------------------------------------------------------
# 7. Extend the sheet with one more trailing blank row (A1022:B1022, same
#    formatting as the row above it)
# ---------------------------------------------------------------------------
$ws.Range("A1021:B1021").Copy()
$ws.Range("A1022:B1022").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 8. Move the active cell / scroll position
# ---------------------------------------------------------------------------
$ws.Range("C32").Select()
